$p = $ppt.ActivePresentation

# Insert a brand-new "Title and Content" slide as the 2nd slide of the
# deck (pushing the former slide 2 "Eszközök:" etc. one position later).
$s = $p.Slides.Add(2, 2)

# Content placeholder gets a short (gibberish) body text.
$content = $s.Shapes.Item(2)
$content.TextFrame.TextRange.Text = "fhihfpeihfűpeihfűpihefűpihpihgűp"
$content.TextFrame.TextRange.LanguageID = "hu-HU"
